# No312. Burst Balloons finished
# Append a new row (row 36) to the leetcode tracker sheet describing the
# newly finished problem, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Write the new row's values -------------------------------------------------
$ws.Range("A36").Value = "No312. Burst Balloons"
$ws.Range("B36").Value = "Hard"
$ws.Range("C36").Value = "https://leetcode.com/problems/burst-balloons/"
$ws.Range("D36").Value = 44539
$ws.Range("E36").Value = "分治"
$ws.Range("F36").Value = "分治的依据是最后一个爆炸的元素为界，因为分治后的结果需要使用分治前的元素，所以需要构造比给定数组左右各多一个空位用于储存爆炸后左右位置的元素"
$ws.Range("G36").Value = "未复习"
$ws.Range("H36").Value = "⭕"

# --- Hyperlink the problem-link cell, same as every other row in C ------------
$ws.Hyperlinks.Add($ws.Range("C36"), "https://leetcode.com/problems/burst-balloons/") | Out-Null

# --- Match formatting of the other "Hard" rows (A:F from row 9, G:H from row 35) -
$ws.Range("A9:F9").Copy() | Out-Null
$ws.Range("A36:F36").PasteSpecial(-4122) | Out-Null
$ws.Range("G35:H35").Copy() | Out-Null
$ws.Range("G36:H36").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row height to fit the wrapped note text, like the neighbouring rows.
$ws.Rows.Item(36).RowHeight = 42

# Leave the selection where the author left it after typing the new row.
$ws.Range("I38").Select() | Out-Null
